# Add a "chr" (chromosome) column to the ratios table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell "chr", styled like the other header cells (bold, bordered, centered)
$headerCell = $ws.Cells.Item(1, 5)
$headerCell.Value = "chr"
$headerCell.Style = $ws.Cells.Item(1, 4).Style

# Chromosome values for each data row, stored as text
$chrValues = @("10", "9", "6", "6", "8", "6", "6", "6", "6", "17", "6", "14")

for ($i = 0; $i -lt $chrValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"       # force text interpretation of numeric-looking strings
    $cell.Value = $chrValues[$i]
    $cell.Style = $ws.Cells.Item($row, 4).Style   # match plain data-cell style, no overrides
}
